$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.Value = "'90.566.11"
$c.ClearFormats()
$ws.Range("E2").Value = '  +1.03%  '

# Row 3
$c = $ws.Range("D3")
$c.Value = "'3.114.09"
$c.ClearFormats()
$ws.Range("E3").Value = '  +2.16%  '

# Row 4
$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.ClearFormats()
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$c = $ws.Range("D5")
$c.Value = "'238.11"
$c.ClearFormats()
$ws.Range("E5").Value = '  +11.68%  '

# Row 6
$c = $ws.Range("D6")
$c.Value = "'625.71"
$c.ClearFormats()
$ws.Range("E6").Value = '  +2.78%  '

# Row 7
$c = $ws.Range("D7")
$c.Value = "'1.12"
$c.ClearFormats()
$ws.Range("E7").Value = '  +5.42%  '

# Row 8
$c = $ws.Range("D8")
$c.Value = "'0.373"
$c.ClearFormats()
$ws.Range("E8").Value = '  +8.41%  '

# Row 9
$ws.Range("E9").Value = '  +0.05%  '

# Row 10
$c = $ws.Range("D10")
$c.Value = "'3.112.35"
$c.ClearFormats()
$ws.Range("E10").Value = '  +2.09%  '

# Row 11
$c = $ws.Range("D11")
$c.Value = "'0.735"
$c.ClearFormats()
$ws.Range("E11").Value = '  +4.13%  '

# Row 12
$c = $ws.Range("D12")
$c.Value = "'0.203"
$c.ClearFormats()
$ws.Range("E12").Value = '  +4.87%  '

# Row 13
$c = $ws.Range("D13")
$c.Value = "'0.0000249"
$c.ClearFormats()
$ws.Range("E13").Value = '  +6.02%  '

# Row 14
$c = $ws.Range("D14")
$c.Value = "'35.28"
$c.ClearFormats()
$ws.Range("E14").Value = '  +4.74%  '

# Row 15
$c = $ws.Range("D15")
$c.Value = "'5.46"
$c.ClearFormats()
$ws.Range("E15").Value = '  +0.34%  '

# Row 16
$c = $ws.Range("D16")
$c.Value = "'90.251.56"
$c.ClearFormats()
$ws.Range("E16").Value = '  +0.84%  '

# Row 17
$c = $ws.Range("D17")
$c.Value = "'3.685.68"
$c.ClearFormats()
$ws.Range("E17").Value = '  +1.33%  '

# Row 18
$c = $ws.Range("D18")
$c.Value = "'3.102.58"
$c.ClearFormats()
$ws.Range("E18").Value = '  +1.44%  '

# Row 19
$c = $ws.Range("D19")
$c.Value = "'3.85"
$c.ClearFormats()
$ws.Range("E19").Value = '  +6.32%  '

# Row 20
$c = $ws.Range("D20")
$c.Value = "'0.0000213"
$c.ClearFormats()
$ws.Range("E20").Value = '  +7.50%  '

# Row 21
$c = $ws.Range("D21")
$c.Value = "'14.26"
$c.ClearFormats()
$ws.Range("E21").Value = '  +2.42%  '

# Row 22
$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D22")
$c.Value = "'5.72"
$c.ClearFormats()
$ws.Range("E22").Value = '  +6.41%  '

# Row 23
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D23")
$c.Value = "'445.57"
$c.ClearFormats()
$ws.Range("E23").Value = '  +0.20%  '

# Row 24
$c = $ws.Range("D24")
$c.Value = "'9.13"
$c.ClearFormats()
$ws.Range("E24").Value = '  +5.12%  '

# Row 25
$c = $ws.Range("D25")
$c.Value = "'5.93"
$c.ClearFormats()
$ws.Range("E25").Value = '  +4.75%  '

# Row 26
$c = $ws.Range("D26")
$c.Value = "'90.40"
$c.ClearFormats()
$ws.Range("E26").Value = '  -0.01%  '

# Row 27
$c = $ws.Range("D27")
$c.Value = "'12.11"
$c.ClearFormats()
$ws.Range("E27").Value = '  +3.62%  '

# Row 28
$c = $ws.Range("D28")
$c.Value = "'3.164.60"
$c.ClearFormats()
$ws.Range("E28").Value = '  -2.94%  '

# Row 29
$ws.Range("E29").Value = '  +0.18%  '

# Row 30
$c = $ws.Range("D30")
$c.Value = "'0.174"
$c.ClearFormats()
$ws.Range("E30").Value = '  +11.33%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D31")
$c.Value = "'0.217"
$c.ClearFormats()
$ws.Range("E31").Value = '  +8.97%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D32")
$c.Value = "'9.27"
$c.ClearFormats()
$ws.Range("E32").Value = '  +2.53%  '

# Row 33
$c = $ws.Range("D33")
$c.Value = "'1.01"
$c.ClearFormats()
$ws.Range("E33").Value = '  +8.72%  '

# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D34")
$c.Value = "'0.107"
$c.ClearFormats()
$ws.Range("E34").Value = '  +27.71%  '

# Row 35
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D35")
$c.Value = "'26.61"
$c.ClearFormats()
$ws.Range("E35").Value = '  -4.18%  '

# Row 36
$ws.Range("B36").Value = 'MantraDAO'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$c = $ws.Range("D36")
$c.Value = "'4.22"
$c.ClearFormats()
$ws.Range("E36").Value = '  +43.11%  '

# Row 37
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D37")
$c.Value = "'0.155"
$c.ClearFormats()
$ws.Range("E37").Value = '  +8.44%  '

# Row 38
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Range("D38")
$c.Value = "'7.31"
$c.ClearFormats()
$ws.Range("E38").Value = '  +9.09%  '

# Row 39
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range("D39")
$c.Value = "'496.26"
$c.ClearFormats()
$ws.Range("E39").Value = '  +1.05%  '

# Row 40
$c = $ws.Range("D40")
$c.Value = "'1.92"
$c.ClearFormats()
$ws.Range("E40").Value = '  +3.03%  '

# Row 41
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D41")
$c.Value = "'3.61"
$c.ClearFormats()
$ws.Range("E41").Value = '  +6.25%  '

# Row 42
$c = $ws.Range("D42")
$c.Value = "'1.29"
$c.ClearFormats()
$ws.Range("E42").Value = '  +2.60%  '

# Row 43
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range("D43")
$c.Value = "'0.417"
$c.ClearFormats()
$ws.Range("E43").Value = '  -0.08%  '

# Row 44
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Range("D44")
$c.Value = "'22.09"
$c.ClearFormats()
$ws.Range("E44").Value = '  -0.24%  '

# Row 45
$ws.Range("E45").Value = '  +0.04%  '

# Row 46
$c = $ws.Range("D46")
$c.Value = "'158.25"
$c.ClearFormats()
$ws.Range("E46").Value = '  +7.36%  '

# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D47")
$c.Value = "'1.89"
$c.ClearFormats()
$ws.Range("E47").Value = '  -0.57%  '

# Row 48
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D48")
$c.Value = "'0.682"
$c.ClearFormats()
$ws.Range("E48").Value = '  +1.11%  '

# Row 49
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D49")
$c.Value = "'4.53"
$c.ClearFormats()
$ws.Range("E49").Value = '  -0.27%  '

# Row 50
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D50")
$c.Value = "'44.85"
$c.ClearFormats()
$ws.Range("E50").Value = '  +0.61%  '

# Row 51
$c = $ws.Range("D51")
$c.Value = "'1.33"
$c.ClearFormats()
$ws.Range("E51").Value = '  +2.29%  '
